$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B held the shared string "R40" (the rule id for the
# last "Good Night" rule). The rule id is being changed to the text "1".
# The leading apostrophe forces Excel to store the digit string as text
# (shared string) instead of re-interpreting it as the number 1.
$ws.Range("B11").Value = "'1"
